$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current values for rows 16-25, columns 1-51 (A-AY)
$data = @{}
for ($r = 16; $r -le 25; $r++) {
    $rowData = @{}
    for ($c = 1; $c -le 51; $c++) {
        $rowData[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $data[$r] = $rowData
}

# Mapping: destination row -> source row (content to move into destination)
$mapping = @{
    16 = 22
    17 = 19
    18 = 20
    19 = 24
    20 = 25
    21 = 16
    22 = 23
    23 = 17
    24 = 18
    25 = 21
}

# Write back permuted rows. Column 9 (I) holds numeric-looking text (e.g. "16"),
# which must stay text, so force text format there before assigning, then clear
# the number-format override so no stray style is left on the cell.
foreach ($dst in 16..25) {
    $src = $mapping[$dst]
    $rowData = $data[$src]
    for ($c = 1; $c -le 51; $c++) {
        $val = $rowData[$c]
        if ($c -eq 9) {
            $cell = $ws.Cells.Item($dst, $c)
            $cell.NumberFormat = "@"
            $cell.Value2 = $val
            $cell.ClearFormats()
        } else {
            $ws.Cells.Item($dst, $c).Value2 = $val
        }
    }
}
